# canasta_diaria.xlsx - "orshanky revised and outliars removed"
# The upstream data-processing pipeline recomputed cantidad_h / cantidad_ajustada /
# population / cal_intake for each "bien" (good) and dropped an outlier item
# ("Hueso de res, pata de res, pata de pollo"). This script reproduces that on the
# live worksheet via COM automation: remove the outlier row, then refresh every
# remaining row (still ordered by descending cal_intake) with its recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the outlier item entirely ---
$outlierName = "Hueso de res, pata de res, pata de pollo"
$usedRows = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $usedRows; $r++) {
    if ($ws.Cells.Item($r, 1).Value2 -eq $outlierName) {
        $ws.Rows.Item($r).Delete()
        break
    }
}

# --- Refreshed values: bien, cantidad_h, cantidad_ajustada, population, cal, cal_intake ---
$refreshed = @(
    @("Harina de maiz", 98.852728790493416, 121.37409532839253, 2976, 383, 464.86279296875),
    @("Arroz, harina de arroz", 93.620393904306553, 114.94969125460553, 2976, 345, 396.57644653320313),
    @("Aceite", 19.696025640934064, 24.183321226648587, 2976, 900, 217.64988708496094),
    @("Azucar", 29.431404478367298, 36.136687085032463, 2976, 393.5, 142.19786071777344),
    @("Queso blanco", 21.704589719413431, 26.649491610065585, 2976, 368.5, 98.203376770019531),
    @("Pastas alimenticias", 46.558181111530594, 57.165414030833908, 2976, 137.5, 78.602447509765625),
    @("Yuca", 33.69095706651288, 41.366682636481457, 2976, 182.33332824707031, 75.425247192382813),
    @("Carne de res (bistec, carne molida, carne para esmechar)", 30.1651314147698, 37.037577105465758, 2976, 196.5, 72.778839111328125),
    @("Lentejas", 17.353111109425946, 21.306626609576647, 2976, 254.5, 54.225364685058594),
    @("Maiz en granos", 10.393337275392266, 12.761224837953685, 2976, 355, 45.302349090576172),
    @("Margarina/Mantequilla", 5.8323734088610575, 7.1611482968093245, 2976, 584, 41.82110595703125),
    @("Cambur", 29.662635196921645, 36.420598387718201, 2976, 113.375, 41.291854858398438),
    @("Carne de pollo", 18.889209435992342, 23.192690318951044, 2976, 174, 40.355281829833984),
    @("Platanos", 17.025249929838282, 20.904069359584522, 2976, 164.85714721679688, 34.46185302734375),
    @("Frijoles", 6.7012290829612366, 8.2279531987764507, 2976, 405.84616088867188, 33.392833709716797),
    @("Leche en polvo, completa o descremada", 6.1035908024798156, 7.4941564862446119, 2976, 428.5, 32.112461090087891),
    @("Huevos (unidades)", 16.378528225806452, 20.110006128588029, 2976, 145, 29.159509658813477),
    @("Pescado fresco", 26.759313247537101, 32.855819927748811, 2976, 85, 27.927446365356445),
    @("Caraotas", 15.459389871166598, 18.981463663039669, 2976, 135.11111450195313, 25.646066665649414),
    @("Papas", 15.005760780906165, 18.424485367472453, 2976, 122.46154022216797, 22.562908172607422),
    @("Cebolla", 13.526305961352522, 16.607970019822481, 2976, 40, 6.6431879997253418),
    @("Pan de trigo", 1.8349654508694526, 2.2530209428520611, 2976, 284.66665649414063, 6.4135994911193848),
    @("Cebollin, ajoporro, cilantro y similares", 6.8738000714971177, 8.439840540770561, 2976, 56.5, 4.7685098648071289),
    @("Aji dulce, pimenton, pimiento", 10.194652644616943, 12.517274501583268, 2976, 35, 4.3810462951660156),
    @("Tomates", 9.1769395345641716, 11.267698496580124, 2976, 23, 2.5915706157684326),
    @("Cafe", 10.527313981325396, 12.925725298623243, 2976, 5, 0.6462862491607666),
    @("Sal", 17.385177094929961, 21.345998065045443, 2976, 0, 0)
)

$r = 2
foreach ($row in $refreshed) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}
